$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1178.4445
$ws.Range("I15").Value = 1178.4445
$ws.Range("K15").Value = 3535.3335
$ws.Range("M15").Value = -3366.3335

$ws.Range("H40").Value = 3350.3
$ws.Range("I40").Value = 2333.6667
$ws.Range("K40").Value = 2333.6667
$ws.Range("M40").Value = -2158.6667

$ws.Range("H43").Value = 1199.5
$ws.Range("J43").Value = 1199.5
$ws.Range("L43").Value = 1199.5
$ws.Range("N43").Value = -1337.5

$ws.Range("H74").Value = 3432.2222
$ws.Range("I74").Value = 1972.5
$ws.Range("K74").Value = 1972.5
$ws.Range("M74").Value = -1036.5

$ws.Range("H77").Value = 3432.2222
$ws.Range("I77").Value = 1972.5
$ws.Range("K77").Value = 9862.5
$ws.Range("M77").Value = -5182.5

$ws.Range("H92").Value = 3977.8147
$ws.Range("I92").Value = 268.52942
$ws.Range("K92").Value = 268.52942
$ws.Range("M92").Value = 979.4705799999999

$ws.Range("H113").Value = 76926920
$ws.Range("I113").Value = 125002744
$ws.Range("K113").Value = 125002744
$ws.Range("M113").Value = -124999490

$ws.Range("H137").Value = 2599.5
$ws.Range("I137").Value = 2244.4443
$ws.Range("K137").Value = 6733.3329
$ws.Range("M137").Value = -4183.3329

$ws.Range("H141").Value = 2377.2727
$ws.Range("I141").Value = 1616.6666
$ws.Range("K141").Value = 4849.9998
$ws.Range("M141").Value = 330.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2538.8572
$ws.Range("I74").Value = 2365.85
$ws.Range("K74").Value = 2365.85
$ws.Range("M74").Value = -1491.85

$ws.Range("H77").Value = 2538.8572
$ws.Range("I77").Value = 2365.85
$ws.Range("K77").Value = 11829.25
$ws.Range("M77").Value = -7461.25

$ws.Range("H97").Value = 986.17645
$ws.Range("I97").Value = 1081.7858
$ws.Range("K97").Value = 1081.7858
$ws.Range("M97").Value = -585.7858000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1653.7
$ws.Range("I105").Value = 1491.6818
$ws.Range("K105").Value = 1491.6818
$ws.Range("M105").Value = 255.3181999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3052.975
$ws.Range("I31").Value = 1417.4375
$ws.Range("J31").Value = 4143.3335
$ws.Range("K31").Value = 1417.4375
$ws.Range("L31").Value = 4143.3335
$ws.Range("M31").Value = -1122.4375
$ws.Range("N31").Value = -4733.3335

$ws.Range("H34").Value = 3052.975
$ws.Range("I34").Value = 1417.4375
$ws.Range("J34").Value = 4143.3335
$ws.Range("K34").Value = 1417.4375
$ws.Range("L34").Value = 4143.3335
$ws.Range("M34").Value = -1215.4375
$ws.Range("N34").Value = -4547.3335

$ws.Range("H58").Value = 2346.5293
$ws.Range("I58").Value = 1291.1538
$ws.Range("K58").Value = 1291.1538
$ws.Range("M58").Value = -1088.1538

$ws.Range("H74").Value = 36367.43
$ws.Range("J74").Value = 36367.43
$ws.Range("L74").Value = 36367.43
$ws.Range("N74").Value = -38115.43

$ws.Range("H77").Value = 36367.43
$ws.Range("J77").Value = 36367.43
$ws.Range("L77").Value = 109102.29
$ws.Range("N77").Value = -117838.29

$ws.Range("H86").Value = 5739.8
$ws.Range("I86").Value = 5889.8
$ws.Range("K86").Value = 5889.8
$ws.Range("M86").Value = -4766.8

$ws.Range("H89").Value = 5739.8
$ws.Range("I89").Value = 5889.8
$ws.Range("K89").Value = 29449
$ws.Range("M89").Value = -23833

$ws.Range("H134").Value = 2782.5293
$ws.Range("I134").Value = 2199.3333
$ws.Range("K134").Value = 6597.999899999999
$ws.Range("M134").Value = -4062.999899999999

$ws.Range("H136").Value = 2346.5293
$ws.Range("I136").Value = 1291.1538
$ws.Range("K136").Value = 3873.4614
$ws.Range("M136").Value = -1323.4614

$ws.Range("H141").Value = 77323.75
$ws.Range("J141").Value = 85512.86
$ws.Range("L141").Value = 85512.86
$ws.Range("N141").Value = -95872.86

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 903.26666
$ws.Range("I5").Value = 577
$ws.Range("J5").Value = 1555.8
$ws.Range("K5").Value = 1731
$ws.Range("L5").Value = 4667.4
$ws.Range("M5").Value = -1619
$ws.Range("N5").Value = -4891.4

$ws.Range("H11").Value = 35983.465
$ws.Range("I11").Value = 54921.58
$ws.Range("J11").Value = 3272.182
$ws.Range("K11").Value = 164764.74
$ws.Range("L11").Value = 9816.545999999998
$ws.Range("M11").Value = -164624.74
$ws.Range("N11").Value = -10096.546

$ws.Range("H132").Value = 2390.25
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 2687
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 24183
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -29243

$ws.Range("H135").Value = 903.26666
$ws.Range("I135").Value = 577
$ws.Range("J135").Value = 1555.8
$ws.Range("K135").Value = 5193
$ws.Range("L135").Value = 14002.2
$ws.Range("M135").Value = -2658
$ws.Range("N135").Value = -19072.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1485.8077
$ws.Range("I97").Value = 1528.7222
$ws.Range("J97").Value = 1389.25
$ws.Range("K97").Value = 1528.7222
$ws.Range("L97").Value = 1389.25
$ws.Range("M97").Value = -1032.7222
$ws.Range("N97").Value = -2381.25

$ws.Range("H102").Value = 2619.1177
$ws.Range("I102").Value = 1359.5
$ws.Range("J102").Value = 3738.7778
$ws.Range("K102").Value = 1359.5
$ws.Range("L102").Value = 3738.7778
$ws.Range("M102").Value = 262.5
$ws.Range("N102").Value = -6982.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 27781394
$ws.Range("I7").Value = 50002828
$ws.Range("J7").Value = 4600.625
$ws.Range("K7").Value = 50002828
$ws.Range("L7").Value = 4600.625
$ws.Range("M7").Value = -50002716
$ws.Range("N7").Value = -4824.625

$ws.Range("H46").Value = 1271.5
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H126").Value = 27781394
$ws.Range("I126").Value = 50002828
$ws.Range("J126").Value = 4600.625
$ws.Range("K126").Value = 150008484
$ws.Range("L126").Value = 13801.875
$ws.Range("M126").Value = -150006014
$ws.Range("N126").Value = -18741.875

$ws.Range("H132").Value = 5857.8
$ws.Range("I132").Value = 3671.3125
$ws.Range("J132").Value = 8356.643
$ws.Range("K132").Value = 11013.9375
$ws.Range("L132").Value = 25069.929
$ws.Range("M132").Value = -8483.9375
$ws.Range("N132").Value = -30129.929

$ws.Range("H136").Value = 2322.95
$ws.Range("I136").Value = 1968.2354
$ws.Range("K136").Value = 5904.706200000001
$ws.Range("M136").Value = -3354.706200000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 80496.5
$ws.Range("J46").Value = 80496.5
$ws.Range("L46").Value = 80496.5
$ws.Range("N46").Value = -80958.5

$ws.Range("H54").Value = 10179.5
$ws.Range("J54").Value = 10179.5
$ws.Range("L54").Value = 10179.5
$ws.Range("N54").Value = -11219.5

$ws.Range("H64").Value = 53777
$ws.Range("I64").Value = 53777
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 53777
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -53529
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 53777
$ws.Range("I67").Value = 53777
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 53777
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -52919
$ws.Range("N67").ClearContents()

$ws.Range("H100").Value = 6961.125
$ws.Range("I100").Value = 7812.857
$ws.Range("K100").Value = 15625.714
$ws.Range("M100").Value = -15084.714

$ws.Range("H132").Value = 4143.6553
$ws.Range("I132").Value = 3868.5652
$ws.Range("K132").Value = 11605.6956
$ws.Range("M132").Value = -9075.695599999999

$ws.Range("H134").Value = 80496.5
$ws.Range("J134").Value = 80496.5
$ws.Range("L134").Value = 241489.5
$ws.Range("N134").Value = -246559.5

$ws.Range("H135").Value = 57402.5
$ws.Range("J135").Value = 57402.5
$ws.Range("L135").Value = 57402.5
$ws.Range("N135").Value = -67542.5

$ws.Range("H136").Value = 3334.2778
$ws.Range("I136").Value = 918.6
$ws.Range("J136").Value = 6353.875
$ws.Range("K136").Value = 2755.8
$ws.Range("L136").Value = 19061.625
$ws.Range("M136").Value = -205.8000000000002
$ws.Range("N136").Value = -24161.625
